# Updates the "展览" (Sheet1) and "全部类型" (Sheet4) worksheets to reflect
# the refreshed show listing: the "常熟·漫魂动漫游戏展01" row is removed,
# every subsequent row shifts up by one, and several "想去人数" (interest
# count) figures are bumped to their newer values. Sheets "演出" and
# "本地生活" are untouched.

function Set-RowData($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i) {
    # Leading apostrophe keeps date-like strings (column B) as plain text
    # instead of letting Excel coerce them into date serial numbers.
    $ws.Cells.Item($row, 2).Value = "'" + $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览": rows 4..26 get the updated data, row 27 is removed.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Set-RowData $ws1 4  "2024.02.16" "太仓·龙狮新春动漫节4.0" "滨河路126号 凯景世纪大酒店" "2024.02.16 08:30-02.16 15:00" 73 45 "https://show.bilibili.com/platform/detail.html?id=81044" "//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg"
Set-RowData $ws1 5  "2024.02.16" "苏州·Good Jump ACG迎新特别篇X动漫品牌博览会" "金山南路288号 广电国际会展中心" "2024.02.16 10:00-02.17 17:00" 12073 60 "https://show.bilibili.com/platform/detail.html?id=79303" "//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg"
Set-RowData $ws1 6  "2024.02.16" "苏州·运动番only专区-Good jump ACG" "金山南路288号 广电国际会展中心" "2024.02.16 10:00-02.17 17:00" 4468 25 "https://show.bilibili.com/platform/detail.html?id=81435" "//i0.hdslb.com/bfs/openplatform/202401/gatL3YjP1706236832019.jpeg"
Set-RowData $ws1 7  "2024.02.24" "张家港·上元节AF 汉服花灯动漫展" "杨舍镇新农路163号 菁英羽毛球馆" "2024.02.24 09:30-02.24 17:00" 40 30 "https://show.bilibili.com/platform/detail.html?id=81508" "//i0.hdslb.com/bfs/openplatform/202401/kd6tKAkO1706524738420.jpeg"
Set-RowData $ws1 8  "2024.02.24" "苏州·世纪幻想动漫游戏展" "相城大道1609号 苏州环球港" "2024.02.24 10:30-02.24 16:30" 59 20 "https://show.bilibili.com/platform/detail.html?id=81565" "//i2.hdslb.com/bfs/openplatform/202401/JFcO48XW1706685060185.jpeg"
Set-RowData $ws1 9  "2024.02.24" "苏州·异境次元 元宵特典 二次元派对狂欢" "盘门路28号 吴门印象20幢 ISZGameZone" "2024.02.24 13:00-02.24 19:00" 30 169 "https://show.bilibili.com/platform/detail.html?id=81807" "//i0.hdslb.com/bfs/openplatform/202402/QcC9Z3MQ1707275914021.jpeg"
Set-RowData $ws1 10 "2024.02.25" "太仓·龙吟动漫游戏展" "滨河路128号 凯景世纪大酒店(太仓滨河路店)" "2024.02.25 10:00-02.25 17:00" 25 45 "https://show.bilibili.com/platform/detail.html?id=81242" "//i2.hdslb.com/bfs/openplatform/202401/lVVrbSra1706508320671.jpeg"
Set-RowData $ws1 11 "2024.02.25" "苏州·第五届次元鹿角动漫游戏展" "清禾路886号 尹山湖大剧院" "2024.02.25 10:00-02.25 17:00" 2586 68 "https://show.bilibili.com/platform/detail.html?id=79333" "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"
Set-RowData $ws1 12 "2024.03.08" "苏州·国风宠物-cosplay展" "金山南路影视城 木渎影视城会展中心" "2024.03.08 09:00-03.10 17:30" 1117 65 "https://show.bilibili.com/platform/detail.html?id=80635" "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"
Set-RowData $ws1 13 "2024.03.17" "苏州·世纪幻想动漫游戏展2.0" "清禾路886号 尹山湖大剧院" "2024.03.17 10:00-03.17 17:00" 181 60 "https://show.bilibili.com/platform/detail.html?id=81387" "//i1.hdslb.com/bfs/openplatform/202401/m0Q2ZB3L1706153205872.jpeg"
Set-RowData $ws1 14 "2024.03.23" "苏州·Look Look动漫嘉年华" "阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店" "2024.03.23 10:00-03.23 17:30" 59 29.9 "https://show.bilibili.com/platform/detail.html?id=81698" "//i0.hdslb.com/bfs/openplatform/202402/ZYkvUFn41706869061984.jpeg"
Set-RowData $ws1 15 "2024.04.04" "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.04.04 10:00-04.05 17:00" 5227 60 "https://show.bilibili.com/platform/detail.html?id=81827" "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"
Set-RowData $ws1 16 "2024.04.06" "苏州·第一届寒假动漫展宅舞比赛-CF01" "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店" "2024.04.06 10:00-04.06 16:00" 63 49 "https://show.bilibili.com/platform/detail.html?id=80528" "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
Set-RowData $ws1 17 "2024.04.13" "苏州·绘时国乙1.0-秩序之外" "石路步行街永福桥浜15号 银河广场" "2024.04.13 13:30-04.13 20:00" 200 78 "https://show.bilibili.com/platform/detail.html?id=80789" "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"
Set-RowData $ws1 18 "2024.04.21" "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.04.21 10:00-04.21 21:00" 543 49.9 "https://show.bilibili.com/platform/detail.html?id=78666" "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"
Set-RowData $ws1 19 "2024.05.01" "昆山·第十二届理想乡动漫游戏展" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.01 10:00-05.03 17:00" 11414 59 "https://show.bilibili.com/platform/detail.html?id=77196" "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"
Set-RowData $ws1 20 "2024.05.01" "苏州·第十七届 I COME ACG  动漫品牌博览会" "金山南路288号 广电国际会展中心" "2024.05.01 10:00-05.02 17:00" 11468 65 "https://show.bilibili.com/platform/detail.html?id=79789" "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"
Set-RowData $ws1 21 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 22 1 "https://show.bilibili.com/platform/detail.html?id=81116" "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"
Set-RowData $ws1 22 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 55 1 "https://show.bilibili.com/platform/detail.html?id=81100" "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"
Set-RowData $ws1 23 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 15 1 "https://show.bilibili.com/platform/detail.html?id=81119" "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"
Set-RowData $ws1 24 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 13 1 "https://show.bilibili.com/platform/detail.html?id=81118" "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"
Set-RowData $ws1 25 "2024.05.03" "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 54 1 "https://show.bilibili.com/platform/detail.html?id=81120" "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"
Set-RowData $ws1 26 "2024.05.03" "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 26 1 "https://show.bilibili.com/platform/detail.html?id=81114" "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

# Drop the now-duplicate last row (row 27); dimension becomes A1:I26.
$ws1.Rows.Item(27).Delete()

# ---------------------------------------------------------------------
# Sheet "全部类型": rows 4..27 get the updated data, row 28 is removed.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Set-RowData $ws4 4  "2024.02.16" "太仓·龙狮新春动漫节4.0" "滨河路126号 凯景世纪大酒店" "2024.02.16 08:30-02.16 15:00" 73 45 "https://show.bilibili.com/platform/detail.html?id=81044" "//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg"
Set-RowData $ws4 5  "2024.02.16" "苏州·Good Jump ACG迎新特别篇X动漫品牌博览会" "金山南路288号 广电国际会展中心" "2024.02.16 10:00-02.17 17:00" 12073 60 "https://show.bilibili.com/platform/detail.html?id=79303" "//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg"
Set-RowData $ws4 6  "2024.02.16" "苏州·运动番only专区-Good jump ACG" "金山南路288号 广电国际会展中心" "2024.02.16 10:00-02.17 17:00" 4468 25 "https://show.bilibili.com/platform/detail.html?id=81435" "//i0.hdslb.com/bfs/openplatform/202401/gatL3YjP1706236832019.jpeg"
Set-RowData $ws4 7  "2024.02.24" "张家港·上元节AF 汉服花灯动漫展" "杨舍镇新农路163号 菁英羽毛球馆" "2024.02.24 09:30-02.24 17:00" 40 30 "https://show.bilibili.com/platform/detail.html?id=81508" "//i0.hdslb.com/bfs/openplatform/202401/kd6tKAkO1706524738420.jpeg"
Set-RowData $ws4 8  "2024.02.24" "苏州·世纪幻想动漫游戏展" "相城大道1609号 苏州环球港" "2024.02.24 10:30-02.24 16:30" 59 20 "https://show.bilibili.com/platform/detail.html?id=81565" "//i2.hdslb.com/bfs/openplatform/202401/JFcO48XW1706685060185.jpeg"
Set-RowData $ws4 9  "2024.02.24" "苏州·异境次元 元宵特典 二次元派对狂欢" "盘门路28号 吴门印象20幢 ISZGameZone" "2024.02.24 13:00-02.24 19:00" 30 169 "https://show.bilibili.com/platform/detail.html?id=81807" "//i0.hdslb.com/bfs/openplatform/202402/QcC9Z3MQ1707275914021.jpeg"
Set-RowData $ws4 10 "2024.02.25" "太仓·龙吟动漫游戏展" "滨河路128号 凯景世纪大酒店(太仓滨河路店)" "2024.02.25 10:00-02.25 17:00" 25 45 "https://show.bilibili.com/platform/detail.html?id=81242" "//i2.hdslb.com/bfs/openplatform/202401/lVVrbSra1706508320671.jpeg"
Set-RowData $ws4 11 "2024.02.25" "苏州·第五届次元鹿角动漫游戏展" "清禾路886号 尹山湖大剧院" "2024.02.25 10:00-02.25 17:00" 2586 68 "https://show.bilibili.com/platform/detail.html?id=79333" "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"
Set-RowData $ws4 12 "2024.03.03" "苏州·龙猫和他的朋友·动漫作品音乐会" "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院" "2024.03.03 19:30-03.03 21:00" 4 60 "https://show.bilibili.com/platform/detail.html?id=81799" "//i0.hdslb.com/bfs/openplatform/202402/gqnOEjvJ1707214629948.jpeg"
Set-RowData $ws4 13 "2024.03.08" "苏州·国风宠物-cosplay展" "金山南路影视城 木渎影视城会展中心" "2024.03.08 09:00-03.10 17:30" 1117 65 "https://show.bilibili.com/platform/detail.html?id=80635" "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"
Set-RowData $ws4 14 "2024.03.17" "苏州·世纪幻想动漫游戏展2.0" "清禾路886号 尹山湖大剧院" "2024.03.17 10:00-03.17 17:00" 181 60 "https://show.bilibili.com/platform/detail.html?id=81387" "//i1.hdslb.com/bfs/openplatform/202401/m0Q2ZB3L1706153205872.jpeg"
Set-RowData $ws4 15 "2024.03.23" "苏州·Look Look动漫嘉年华" "阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店" "2024.03.23 10:00-03.23 17:30" 59 29.9 "https://show.bilibili.com/platform/detail.html?id=81698" "//i0.hdslb.com/bfs/openplatform/202402/ZYkvUFn41706869061984.jpeg"
Set-RowData $ws4 16 "2024.04.04" "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会" "金山南路288号木渎影视城F2 苏州广电国际会展中心" "2024.04.04 10:00-04.05 17:00" 5227 60 "https://show.bilibili.com/platform/detail.html?id=81827" "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg"
Set-RowData $ws4 17 "2024.04.06" "苏州·第一届寒假动漫展宅舞比赛-CF01" "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店" "2024.04.06 10:00-04.06 16:00" 63 49 "https://show.bilibili.com/platform/detail.html?id=80528" "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"
Set-RowData $ws4 18 "2024.04.13" "苏州·绘时国乙1.0-秩序之外" "石路步行街永福桥浜15号 银河广场" "2024.04.13 13:30-04.13 20:00" 200 78 "https://show.bilibili.com/platform/detail.html?id=80789" "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"
Set-RowData $ws4 19 "2024.04.21" "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0" "清禾路888号2号楼3楼 格莱美婚礼宴会中心" "2024.04.21 10:00-04.21 21:00" 543 49.9 "https://show.bilibili.com/platform/detail.html?id=78666" "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"
Set-RowData $ws4 20 "2024.05.01" "昆山·第十二届理想乡动漫游戏展" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.01 10:00-05.03 17:00" 11414 59 "https://show.bilibili.com/platform/detail.html?id=77196" "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"
Set-RowData $ws4 21 "2024.05.01" "苏州·第十七届 I COME ACG  动漫品牌博览会" "金山南路288号 广电国际会展中心" "2024.05.01 10:00-05.02 17:00" 11468 65 "https://show.bilibili.com/platform/detail.html?id=79789" "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"
Set-RowData $ws4 22 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 22 1 "https://show.bilibili.com/platform/detail.html?id=81116" "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg"
Set-RowData $ws4 23 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 55 1 "https://show.bilibili.com/platform/detail.html?id=81100" "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg"
Set-RowData $ws4 24 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 15 1 "https://show.bilibili.com/platform/detail.html?id=81119" "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg"
Set-RowData $ws4 25 "2024.05.02" "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.02 14:00-05.02 16:00" 13 1 "https://show.bilibili.com/platform/detail.html?id=81118" "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg"
Set-RowData $ws4 26 "2024.05.03" "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 54 1 "https://show.bilibili.com/platform/detail.html?id=81120" "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg"
Set-RowData $ws4 27 "2024.05.03" "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会" "花桥经济开发区绿地大道1598号 花桥国际博览中心" "2024.05.03 14:00-05.03 16:00" 26 1 "https://show.bilibili.com/platform/detail.html?id=81114" "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg"

# Drop the now-duplicate last row (row 28); dimension becomes A1:I27.
$ws4.Rows.Item(28).Delete()
